$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    479.3767310549769,
    473.5485038290082,
    469.2932929477648,
    633.3911647045371,
    434.9692696308142,
    224.4069150760089,
    590.8390558921017,
    4594.19590353183,
    4343.158628414675,
    5290.860642359125,
    5300.480242313406,
    906.8439227340259,
    4260.484519373642,
    871.0477174534225,
    1048.526311578465,
    1046.549957709452
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
